# Check in for December 28 Task
# Add a new "Resource" column (with header + sample value) to the
# API-Testing sheet, matching the formatting already used by the
# neighbouring header/data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Value = "Resource"
$ws.Range("Q2").Value = "Test"

# Match Q1's formatting to the rest of the header row (column P's header
# style) the way a user would via copy / paste-special-formats.
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
